$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3600
$ws.Range("I4").Value = 3600
$ws.Range("K4").Value = 3600
$ws.Range("M4").Value = -3486
$ws.Range("H9").Value = 471.66666
$ws.Range("I9").Value = 471.66666
$ws.Range("K9").Value = 471.66666
$ws.Range("M9").Value = -302.66666
$ws.Range("H15").Value = 106.375
$ws.Range("I15").Value = 106.375
$ws.Range("K15").Value = 319.125
$ws.Range("M15").Value = -150.125
$ws.Range("H32").Value = 9836.125
$ws.Range("J32").Value = 9916.666999999999
$ws.Range("L32").Value = 9916.666999999999
$ws.Range("N32").Value = -10568.667
$ws.Range("H55").Value = 409.6
$ws.Range("J55").Value = 478.5
$ws.Range("L55").Value = 478.5
$ws.Range("N55").Value = -906.5
$ws.Range("H92").Value = 1249.2307
$ws.Range("I92").Value = 1395.8334
$ws.Range("J92").Value = 1123.5714
$ws.Range("K92").Value = 1395.8334
$ws.Range("L92").Value = 1123.5714
$ws.Range("M92").Value = -147.8334
$ws.Range("N92").Value = -3619.5714
$ws.Range("H99").Value = 402.8
$ws.Range("I99").Value = 402.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1208.4
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 289.5999999999999
$ws.Range("N99").ClearContents()
$ws.Range("H112").Value = 4216
$ws.Range("J112").Value = 4216
$ws.Range("L112").Value = 12648
$ws.Range("N112").Value = -14864
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H125").Value = 1500
$ws.Range("J125").Value = 1500
$ws.Range("L125").Value = 13500
$ws.Range("N125").Value = -18420
$ws.Range("H129").Value = 1192
$ws.Range("I129").Value = 1057.3334
$ws.Range("K129").Value = 3172.0002
$ws.Range("M129").Value = 1827.9998
$ws.Range("H135").Value = 1768.4445
$ws.Range("I135").Value = 1645.1428
$ws.Range("K135").Value = 14806.2852
$ws.Range("M135").Value = -12271.2852
$ws.Range("H137").Value = 2111
$ws.Range("J137").Value = 1950
$ws.Range("L137").Value = 5850
$ws.Range("N137").Value = -10950

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8000
$ws.Range("I2").Value = 8000
$ws.Range("K2").Value = 8000
$ws.Range("M2").Value = -7887
$ws.Range("H74").Value = 5200
$ws.Range("I74").Value = 2750
$ws.Range("K74").Value = 2750
$ws.Range("M74").Value = -1876
$ws.Range("H77").Value = 5200
$ws.Range("I77").Value = 2750
$ws.Range("K77").Value = 13750
$ws.Range("M77").Value = -9382
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H116").Value = 8000
$ws.Range("I116").Value = 8000
$ws.Range("K116").Value = 8000
$ws.Range("M116").Value = -5706

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8000
$ws.Range("I3").Value = 8000
$ws.Range("K3").Value = 8000
$ws.Range("M3").Value = -7886
$ws.Range("H107").Value = 4811
$ws.Range("I107").Value = 4811
$ws.Range("K107").Value = 4811
$ws.Range("M107").Value = -2891

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1466.6666
$ws.Range("I16").Value = 1200
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -913
$ws.Range("N16").Value = -2574
$ws.Range("H31").Value = 4891.125
$ws.Range("I31").Value = 4464
$ws.Range("K31").Value = 4464
$ws.Range("M31").Value = -4169
$ws.Range("H34").Value = 4891.125
$ws.Range("I34").Value = 4464
$ws.Range("K34").Value = 4464
$ws.Range("M34").Value = -4262
$ws.Range("H102").Value = 44666.668
$ws.Range("J102").Value = 44666.668
$ws.Range("L102").Value = 44666.668
$ws.Range("N102").Value = -49534.668
$ws.Range("H107").Value = 767
$ws.Range("I107").Value = 900.5
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 900.5
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1019.5
$ws.Range("N107").Value = -4340
$ws.Range("H113").Value = 1466.6666
$ws.Range("I113").Value = 1200
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -6340

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1154077.5
$ws.Range("I4").Value = 1363809.9
$ws.Range("K4").Value = 4091429.7
$ws.Range("M4").Value = -4091317.7
$ws.Range("H12").Value = 766.6667
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 766.6667
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 2300.0001
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -2646.0001
$ws.Range("H131").Value = 3709.875
$ws.Range("I131").Value = 1230
$ws.Range("K131").Value = 3690
$ws.Range("M131").Value = 1350

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5500
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 5500
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H107").Value = 483.66666
$ws.Range("I107").Value = 224.75
$ws.Range("J107").Value = 1001.5
$ws.Range("K107").Value = 224.75
$ws.Range("L107").Value = 1001.5
$ws.Range("M107").Value = 1695.25
$ws.Range("N107").Value = -4841.5
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102080
$ws.Range("H113").Value = 8500
$ws.Range("I113").Value = 9000
$ws.Range("K113").Value = 9000
$ws.Range("M113").Value = -6830

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 19333.334
$ws.Range("J99").Value = 18000
$ws.Range("L99").Value = 18000
$ws.Range("N99").Value = -23990
$ws.Range("H136").Value = 8700
$ws.Range("I136").Value = 8700
$ws.Range("K136").Value = 26100
$ws.Range("M136").Value = -23550

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H30").Value = 9009
$ws.Range("I30").Value = 9009
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 9009
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -8902
$ws.Range("N30").ClearContents()
$ws.Range("H100").Value = 2646.6155
$ws.Range("I100").Value = 3334.25
$ws.Range("J100").Value = 1546.4
$ws.Range("K100").Value = 6668.5
$ws.Range("L100").Value = 3092.8
$ws.Range("M100").Value = -6127.5
$ws.Range("N100").Value = -4174.8
$ws.Range("H107").Value = 2955.5
$ws.Range("I107").Value = 1607.3334
$ws.Range("K107").Value = 4822.0002
$ws.Range("M107").Value = -2902.0002
$ws.Range("H136").Value = 789
$ws.Range("I136").Value = 786.25
$ws.Range("K136").Value = 2358.75
$ws.Range("M136").Value = 191.25

Write-Host "All edits applied"